# Auto-generated edit script: updates Leve price/profit figures across all sheets
# per the scheduled market-data refresh (current prices, Leve profit recalculation).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow | Beeswax
$ws.Range("H12").Value = 418.83334
$ws.Range("I12").Value = 364.1111
$ws.Range("K12").Value = 364.1111
$ws.Range("M12").Value = -194.1111

# Row 15: Morning Glass of Ether | Ether
$ws.Range("H15").Value = 2420.3264
$ws.Range("I15").Value = 2420.3264
$ws.Range("K15").Value = 7260.9792
$ws.Range("M15").Value = -7091.9792

# Row 17: One for the Road | Potion
$ws.Range("H17").Value = 1627012.1
$ws.Range("J17").Value = 1627012.1
$ws.Range("L17").Value = 4881036.300000001
$ws.Range("N17").Value = -4881372.300000001

# Row 32: Automata for the People | Crab Oil
$ws.Range("H32").Value = 4755
$ws.Range("J32").Value = 4755
$ws.Range("L32").Value = 4755
$ws.Range("N32").Value = -5407

# Row 42: Eye of the Beholder | Hi-Potion of Dexterity
$ws.Range("H42").Value = 156.27272
$ws.Range("I42").Value = 24.333334
$ws.Range("K42").Value = 73.00000199999999
$ws.Range("M42").Value = 156.999998

# Row 47: Open Your Grimoire to Page 42 | Embossed Book of Silver
$ws.Range("H47").Value = 8795
$ws.Range("I47").Value = 4855.6665
$ws.Range("J47").Value = 11749.5
$ws.Range("K47").Value = 4855.6665
$ws.Range("L47").Value = 11749.5
$ws.Range("M47").Value = -3883.6665
$ws.Range("N47").Value = -13693.5

# Row 74: Adhesive of Antipathy | Wing Glue
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

# Row 77: It's Gonna Grow Back (L) | Wing Glue
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 28613.908
$ws.Range("I132").Value = 23305.223
$ws.Range("J132").Value = 52503
$ws.Range("K132").Value = 69915.66900000001
$ws.Range("L132").Value = 157509
$ws.Range("M132").Value = -67385.66900000001
$ws.Range("N132").Value = -162569

# Row 133: Big Brush, Big Dreams | Ginseng Angle Brush
$ws.Range("H133").Value = 63964.168
$ws.Range("J133").Value = 63964.168
$ws.Range("L133").Value = 63964.168
$ws.Range("N133").Value = -74084.16800000001

# Row 140: Tome for Tradition | Book of Ra'Kaznar
$ws.Range("H140").Value = 90593
$ws.Range("J140").Value = 90593
$ws.Range("L140").Value = 90593
$ws.Range("N140").Value = -100953

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Range("H2").Value = 4738.3335
$ws.Range("I2").Value = 840
$ws.Range("K2").Value = 840
$ws.Range("M2").Value = -727

# Row 32: Ingot We Trust | Steel Ingot
$ws.Range("H32").Value = 9851.691000000001
$ws.Range("I32").Value = 3487.6
$ws.Range("J32").Value = 22303.174
$ws.Range("K32").Value = 3487.6
$ws.Range("L32").Value = 22303.174
$ws.Range("M32").Value = -3200.6
$ws.Range("N32").Value = -22877.174

# Row 44: Very Slow Array | Mythril Plate
$ws.Range("H44").Value = 46000
$ws.Range("I44").Value = 46000
$ws.Range("K44").Value = 46000
$ws.Range("M44").Value = -45512

# Row 55: Employee Retention | Mythril Elmo
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

# Row 113: Catching an Earful | Bluespirit Headgear of Casting
$ws.Range("H113").Value = 90000
$ws.Range("J113").Value = 90000
$ws.Range("L113").Value = 90000
$ws.Range("N113").Value = -98678

# Row 116: No Scope | Titanbronze Ingot
$ws.Range("H116").Value = 4738.3335
$ws.Range("I116").Value = 840
$ws.Range("K116").Value = 840
$ws.Range("M116").Value = 1454

# Row 130: A Gift of Gloves | Chondrite Gloves of Casting
$ws.Range("H130").Value = 53927.125
$ws.Range("J130").Value = 53927.125
$ws.Range("L130").Value = 53927.125
$ws.Range("N130").Value = -63967.125

# Row 141: Essays on Equipment | Ra'Kaznar Greaves of Maiming
$ws.Range("H141").Value = 35390
$ws.Range("I141").Value = 35390
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 35390
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -30210

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells | Bronze Ingot
$ws.Range("H3").Value = 4738.3335
$ws.Range("I3").Value = 840
$ws.Range("K3").Value = 840
$ws.Range("M3").Value = -726

# Row 109: Here Comes the Hammer | Deepgold Sledgehammer
$ws.Range("H109").Value = 49994.75
$ws.Range("J109").Value = 49994.75
$ws.Range("L109").Value = 49994.75
$ws.Range("N109").Value = -52768.75

# Row 137: Dagger Swagger | Cobalt Tungsten Khukuri
$ws.Range("H137").Value = 120000
$ws.Range("J137").Value = 120000
$ws.Range("L137").Value = 120000
$ws.Range("N137").Value = -130200

$ws = $wb.Worksheets.Item("CRP")
# Row 48: The Cold, Cold Ground | Oak Pattens
$ws.Range("H48").Value = 54561.332
$ws.Range("J48").Value = 54561.332
$ws.Range("L48").Value = 54561.332
$ws.Range("N48").Value = -55513.332

# Row 70: A Reward Fitting of the Faithful | Holy Cedar Necklace
$ws.Range("H70").Value = 17000
$ws.Range("J70").Value = 17000
$ws.Range("L70").Value = 17000
$ws.Range("N70").Value = -17630

# Row 73: Just Rewards for Just Devotion (L) | Holy Cedar Necklace
$ws.Range("H73").Value = 17000
$ws.Range("J73").Value = 17000
$ws.Range("L73").Value = 17000
$ws.Range("N73").Value = -19184

# Row 86: Birch, Please | Birch Lumber
$ws.Range("H86").Value = 12832.5
$ws.Range("I86").Value = 17980.75
$ws.Range("J86").Value = 7684.25
$ws.Range("K86").Value = 17980.75
$ws.Range("L86").Value = 7684.25
$ws.Range("M86").Value = -16857.75
$ws.Range("N86").Value = -9930.25

# Row 89: Built This City on Blocks and Soul (L) | Birch Lumber
$ws.Range("H89").Value = 12832.5
$ws.Range("I89").Value = 17980.75
$ws.Range("J89").Value = 7684.25
$ws.Range("K89").Value = 89903.75
$ws.Range("L89").Value = 38421.25
$ws.Range("M89").Value = -84287.75
$ws.Range("N89").Value = -49653.25

# Row 94: Beech, Please | Beech Lumber
$ws.Range("H94").Value = 60012
$ws.Range("I94").Value = 60012
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 60012
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -59561
$ws.Range("N94").ClearContents()

# Row 105: Zelkova, My Love | Zelkova Lumber
$ws.Range("H105").Value = 20266.125
$ws.Range("I105").Value = 34710.332
$ws.Range("K105").Value = 34710.332
$ws.Range("M105").Value = -32963.332

# Row 107: Built to Last | White Oak Lumber
$ws.Range("H107").Value = 3852.7144
$ws.Range("I107").Value = 1842.9166
$ws.Range("J107").Value = 6532.4443
$ws.Range("K107").Value = 1842.9166
$ws.Range("L107").Value = 6532.4443
$ws.Range("M107").Value = 77.08339999999998
$ws.Range("N107").Value = -10372.4443

# Row 122: Timber of Tenkonto | Horse Chestnut Lumber
$ws.Range("H122").Value = 6209.857
$ws.Range("I122").Value = 3126.6667
$ws.Range("K122").Value = 9380.000100000001
$ws.Range("M122").Value = -6930.000100000001

# Row 132: Hull Lotta Damage | Ginseng Lumber
$ws.Range("H132").Value = 9537.611000000001
$ws.Range("I132").Value = 2560.2727
$ws.Range("K132").Value = 7680.8181
$ws.Range("M132").Value = -5150.8181

# Row 141: No Greater Treasure | Claro Walnut Necklace of Gathering
$ws.Range("H141").Value = 253404.23
$ws.Range("J141").Value = 253404.23
$ws.Range("L141").Value = 253404.23
$ws.Range("N141").Value = -263764.23

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food | Table Salt
$ws.Range("H2").Value = 87.60869599999999
$ws.Range("I2").Value = 92.64706
$ws.Range("J2").Value = 73.333336
$ws.Range("K2").Value = 555.8823599999999
$ws.Range("L2").Value = 440.000016
$ws.Range("M2").Value = -442.8823599999999
$ws.Range("N2").Value = -666.000016

# Row 38: Pretty as a Picture | Dark Vinegar
$ws.Range("H38").Value = 48.625
$ws.Range("I38").Value = 11.636364
$ws.Range("J38").Value = 130
$ws.Range("K38").Value = 34.909092
$ws.Range("L38").Value = 390
$ws.Range("M38").Value = 312.090908
$ws.Range("N38").Value = -1084

# Row 113: Can't Eat Just One | Night Vinegar
$ws.Range("H113").Value = 8627.058999999999
$ws.Range("I113").Value = 15964.625
$ws.Range("J113").Value = 2104.7778
$ws.Range("K113").Value = 47893.875
$ws.Range("L113").Value = 6314.3334
$ws.Range("M113").Value = -45723.875
$ws.Range("N113").Value = -10654.3334

# Row 115: Mixology | Blood Tomato Juice
$ws.Range("H115").Value = 2566.6667
$ws.Range("I115").Value = 2200
$ws.Range("K115").Value = 6600
$ws.Range("M115").Value = -5425

# Row 122: Salt of the North | Northern Sea Salt
$ws.Range("H122").Value = 8279799.5
$ws.Range("I122").Value = 13348777
$ws.Range("K122").Value = 120138993
$ws.Range("M122").Value = -120136543

# Row 124: Bobbing for Compliments | Island Miq'abob
$ws.Range("H124").Value = 3362.3333
$ws.Range("I124").Value = 2834.8
$ws.Range("K124").Value = 8504.400000000001
$ws.Range("M124").Value = -3594.400000000001

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 1490.9
$ws.Range("I131").Value = 1319.25
$ws.Range("J131").Value = 1498.0521
$ws.Range("K131").Value = 3957.75
$ws.Range("L131").Value = 4494.156300000001
$ws.Range("M131").Value = 1082.25
$ws.Range("N131").Value = -14574.1563

# Row 132: More Mezcal | Cooking Mezcal
$ws.Range("H132").Value = 2166.3333
$ws.Range("I132").Value = 1956.8572
$ws.Range("K132").Value = 17611.7148
$ws.Range("M132").Value = -15081.7148

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 21641.867
$ws.Range("I70").Value = 29429
$ws.Range("K70").Value = 29429
$ws.Range("M70").Value = -29159

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 21641.867
$ws.Range("I73").Value = 29429
$ws.Range("K73").Value = 29429
$ws.Range("M73").Value = -28493

# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Range("H97").Value = 6764.3887
$ws.Range("I97").Value = 1697.5834
$ws.Range("J97").Value = 16898
$ws.Range("K97").Value = 1697.5834
$ws.Range("L97").Value = 16898
$ws.Range("M97").Value = -1201.5834
$ws.Range("N97").Value = -17890

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Range("H126").Value = 12709.272
$ws.Range("I126").Value = 17758.428
$ws.Range("J126").Value = 10353
$ws.Range("K126").Value = 53275.284
$ws.Range("L126").Value = 31059
$ws.Range("M126").Value = -50805.284
$ws.Range("N126").Value = -35999

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic | Boar Leather
$ws.Range("H46").Value = 2564.3
$ws.Range("I46").Value = 1609.4
$ws.Range("J46").Value = 3519.2
$ws.Range("K46").Value = 1609.4
$ws.Range("L46").Value = 3519.2
$ws.Range("M46").Value = -1421.4
$ws.Range("N46").Value = -3895.2

# Row 112: A Slippery Slope | Gliderskin Boots of Casting
$ws.Range("H112").Value = 80791.336
$ws.Range("J112").Value = 80791.336
$ws.Range("L112").Value = 80791.336
$ws.Range("N112").Value = -83745.336

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Range("H132").Value = 2879515.2
$ws.Range("I132").Value = 11600.571
$ws.Range("K132").Value = 34801.713
$ws.Range("M132").Value = -32271.713

$ws = $wb.Worksheets.Item("WVR")
# Row 113: A Tender Table | Pixie Floss
$ws.Range("H113").Value = 5217.5386
$ws.Range("I113").Value = 6647.6665
$ws.Range("K113").Value = 19942.9995
$ws.Range("M113").Value = -17772.9995

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 6443
$ws.Range("I122").Value = 2634.7693
$ws.Range("K122").Value = 7904.3079
$ws.Range("M122").Value = -5454.3079

# Row 124: Hot Heads | Almasty Serge Hat of Casting
$ws.Range("H124").Value = 48875
$ws.Range("J124").Value = 48875
$ws.Range("L124").Value = 48875
$ws.Range("N124").Value = -58695
